$wb = $excel.ActiveWorkbook

# --- Sheet "kite" (sheet1): nomenclature cleanup ---
$kite = $wb.Worksheets.Item("kite")

# Row4 ("obgen.p" = 120) is renamed to "obGen.p" (capitalisation fix).
$kite.Range("A4").Value = "obGen.p"

# Row5 ("prop.p" = 120) is redundant with the renamed row4 and is removed,
# shifting the rows below (obBatt.p, avio.C, and the trailing blank
# formatted cell) up by one.
$kite.Rows("5").Delete()

# Restore the cursor position recorded in the saved file (kept pointing at
# the row just below the last data row, even though it's now out of the
# used range).
$kite.Range("B14").Select()

# --- Sheet "tether" (sheet2): rename "sigma" -> "sigma_max" ---
$tether = $wb.Worksheets.Item("tether")
$tether.Range("A9").Value = "sigma_max"

# Restore saved selection / make this the active tab.
$tether.Range("C14").Select()
$tether.Activate()
